$d = $word.ActiveDocument

# Replace the Jinja variable name "messages_distribution" with "image"
# inside the {{ ... }} placeholder (keeps it a single run for now).
$d.Content.Find.Execute("messages_distribution", $true, $false, $false, $false, $false,
                         $true, 1, $false, "image", 2)

# Force the run containing "image" to split off from the surrounding
# "{{ " / " }}" text by nudging (and reverting) a character formatting
# property on just that sub-range. The engine keeps the resulting run
# boundaries even once the toggled property is restored, giving us three
# runs: "{{ ", "image", " }}" - matching the decorator-ready markup the
# image-rendering template expects.
$rng = $d.Content.Duplicate
$rng.Find.Execute("image", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Bold = 1
$rng.Font.Bold = 0
